# Updates market-price-derived columns (H..N) on several Leve rows
# across all 8 sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(46, 8).Value = 7999  # ALC!H46 was 0
$ws.Cells.Item(46, 10).Value = 7999  # ALC!J46 was 0
$ws.Cells.Item(46, 12).Value = 23997  # ALC!L46 was 0
$ws.Cells.Item(46, 14).Value = -24235  # ALC!N46 was None

$ws.Cells.Item(60, 8).Value = 7999  # ALC!H60 was 0
$ws.Cells.Item(60, 10).Value = 7999  # ALC!J60 was 0
$ws.Cells.Item(60, 12).Value = 23997  # ALC!L60 was 0
$ws.Cells.Item(60, 14).Value = -24965  # ALC!N60 was None

$ws.Cells.Item(74, 8).Value = 8889  # ALC!H74 was 8664
$ws.Cells.Item(74, 9).Value = 7500.5  # ALC!I74 was 6658.6665
$ws.Cells.Item(74, 10).Value = 9285.714  # ALC!J74 was 9666.666999999999
$ws.Cells.Item(74, 11).Value = 7500.5  # ALC!K74 was 6658.6665
$ws.Cells.Item(74, 12).Value = 9285.714  # ALC!L74 was 9666.666999999999
$ws.Cells.Item(74, 13).Value = -6564.5  # ALC!M74 was -5722.6665
$ws.Cells.Item(74, 14).Value = -11157.714  # ALC!N74 was -11538.667

$ws.Cells.Item(77, 8).Value = 8889  # ALC!H77 was 8664
$ws.Cells.Item(77, 9).Value = 7500.5  # ALC!I77 was 6658.6665
$ws.Cells.Item(77, 10).Value = 9285.714  # ALC!J77 was 9666.666999999999
$ws.Cells.Item(77, 11).Value = 37502.5  # ALC!K77 was 33293.3325
$ws.Cells.Item(77, 12).Value = 46428.57  # ALC!L77 was 48333.335
$ws.Cells.Item(77, 13).Value = -32822.5  # ALC!M77 was -28613.3325
$ws.Cells.Item(77, 14).Value = -55788.57  # ALC!N77 was -57693.335

$ws.Cells.Item(135, 8).Value = 2076.44  # ALC!H135 was 1946.6666
$ws.Cells.Item(135, 10).Value = 3468.5  # ALC!J135 was 3823.875
$ws.Cells.Item(135, 12).Value = 31216.5  # ALC!L135 was 34414.875
$ws.Cells.Item(135, 14).Value = -36286.5  # ALC!N135 was -39484.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2078.1428  # ARM!H45 was 2130.625
$ws.Cells.Item(45, 10).Value = 2499.3333  # ARM!J45 was 2499
$ws.Cells.Item(45, 12).Value = 2499.3333  # ARM!L45 was 2499
$ws.Cells.Item(45, 14).Value = -3253.3333  # ARM!N45 was -3253

$ws.Cells.Item(97, 8).Value = 1486241.4  # ARM!H97 was 1688801.5
$ws.Cells.Item(97, 9).Value = 2475064.2  # ARM!I97 was 2855712.8
$ws.Cells.Item(97, 10).Value = 3006.8  # ARM!J97 was 3263.111
$ws.Cells.Item(97, 11).Value = 2475064.2  # ARM!K97 was 2855712.8
$ws.Cells.Item(97, 12).Value = 3006.8  # ARM!L97 was 3263.111
$ws.Cells.Item(97, 13).Value = -2474568.2  # ARM!M97 was -2855216.8
$ws.Cells.Item(97, 14).Value = -3998.8  # ARM!N97 was -4255.111

$ws.Cells.Item(102, 8).Value = 11954483  # ARM!H102 was 12552157
$ws.Cells.Item(102, 9).Value = 2396.4443  # ARM!I102 was 2478.647
$ws.Cells.Item(102, 11).Value = 2396.4443  # ARM!K102 was 2478.647
$ws.Cells.Item(102, 13).Value = -774.4443000000001  # ARM!M102 was -856.6469999999999

$ws.Cells.Item(122, 8).Value = 4123.5835  # ARM!H122 was 4549.2
$ws.Cells.Item(122, 9).Value = 2622.75  # ARM!I122 was 3250
$ws.Cells.Item(122, 11).Value = 7868.25  # ARM!K122 was 9750
$ws.Cells.Item(122, 13).Value = -5418.25  # ARM!M122 was -7300

$ws.Cells.Item(132, 8).Value = 3914.723  # ARM!H132 was 3870.5
$ws.Cells.Item(132, 9).Value = 2834.8955  # ARM!I132 was 2796.147
$ws.Cells.Item(132, 11).Value = 8504.6865  # ARM!K132 was 8388.440999999999
$ws.Cells.Item(132, 13).Value = -5974.6865  # ARM!M132 was -5858.440999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44, 8).Value = 0  # BSM!H44 was 25000
$ws.Cells.Item(44, 10).Value = 0  # BSM!J44 was 25000
$ws.Cells.Item(44, 12).Value = 0  # BSM!L44 was 25000
$ws.Cells.Item(44, 14).ClearContents()  # BSM!N44 was -25994

$ws.Cells.Item(94, 8).Value = 1490.2858  # BSM!H94 was 1532.2222
$ws.Cells.Item(94, 10).Value = 1544.2727  # BSM!J94 was 1662.9
$ws.Cells.Item(94, 12).Value = 1544.2727  # BSM!L94 was 1662.9
$ws.Cells.Item(94, 14).Value = -2446.2727  # BSM!N94 was -2564.9

$ws.Cells.Item(107, 8).Value = 1362  # BSM!H107 was 1395.5
$ws.Cells.Item(107, 9).Value = 1288.125  # BSM!I107 was 1320.6428
$ws.Cells.Item(107, 11).Value = 1288.125  # BSM!K107 was 1320.6428
$ws.Cells.Item(107, 13).Value = 631.875  # BSM!M107 was 599.3571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1107.9166  # CRP!H22 was 903.86664
$ws.Cells.Item(22, 9).Value = 252.11111  # CRP!I22 was 277.375
$ws.Cells.Item(22, 10).Value = 3675.3333  # CRP!J22 was 1619.8572
$ws.Cells.Item(22, 11).Value = 252.11111  # CRP!K22 was 277.375
$ws.Cells.Item(22, 12).Value = 3675.3333  # CRP!L22 was 1619.8572
$ws.Cells.Item(22, 13).Value = 97.88889  # CRP!M22 was 72.625
$ws.Cells.Item(22, 14).Value = -4375.3333  # CRP!N22 was -2319.8572

$ws.Cells.Item(31, 8).Value = 40004320  # CRP!H31 was 37041176
$ws.Cells.Item(31, 10).Value = 5577.5  # CRP!J31 was 5209.75
$ws.Cells.Item(31, 12).Value = 5577.5  # CRP!L31 was 5209.75
$ws.Cells.Item(31, 14).Value = -6167.5  # CRP!N31 was -5799.75

$ws.Cells.Item(34, 8).Value = 40004320  # CRP!H34 was 37041176
$ws.Cells.Item(34, 10).Value = 5577.5  # CRP!J34 was 5209.75
$ws.Cells.Item(34, 12).Value = 5577.5  # CRP!L34 was 5209.75
$ws.Cells.Item(34, 14).Value = -5981.5  # CRP!N34 was -5613.75

$ws.Cells.Item(58, 8).Value = 4431.7417  # CRP!H58 was 4415.75
$ws.Cells.Item(58, 9).Value = 3488.8125  # CRP!I58 was 3425.2778
$ws.Cells.Item(58, 10).Value = 5437.533  # CRP!J58 was 5689.2144
$ws.Cells.Item(58, 11).Value = 3488.8125  # CRP!K58 was 3425.2778
$ws.Cells.Item(58, 12).Value = 5437.533  # CRP!L58 was 5689.2144
$ws.Cells.Item(58, 13).Value = -3285.8125  # CRP!M58 was -3222.2778
$ws.Cells.Item(58, 14).Value = -5843.533  # CRP!N58 was -6095.2144

$ws.Cells.Item(107, 8).Value = 925  # CRP!H107 was 683.3333
$ws.Cells.Item(107, 9).Value = 925  # CRP!I107 was 668.75
$ws.Cells.Item(107, 10).Value = 0  # CRP!J107 was 800
$ws.Cells.Item(107, 11).Value = 925  # CRP!K107 was 668.75
$ws.Cells.Item(107, 12).Value = 0  # CRP!L107 was 800
$ws.Cells.Item(107, 13).Value = 995  # CRP!M107 was 1251.25
$ws.Cells.Item(107, 14).ClearContents()  # CRP!N107 was -4640

$ws.Cells.Item(136, 8).Value = 4431.7417  # CRP!H136 was 4415.75
$ws.Cells.Item(136, 9).Value = 3488.8125  # CRP!I136 was 3425.2778
$ws.Cells.Item(136, 10).Value = 5437.533  # CRP!J136 was 5689.2144
$ws.Cells.Item(136, 11).Value = 10466.4375  # CRP!K136 was 10275.8334
$ws.Cells.Item(136, 12).Value = 16312.599  # CRP!L136 was 17067.6432
$ws.Cells.Item(136, 13).Value = -7916.4375  # CRP!M136 was -7725.8334
$ws.Cells.Item(136, 14).Value = -21412.599  # CRP!N136 was -22167.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 1465.9166  # CUL!H107 was 1335
$ws.Cells.Item(107, 10).Value = 1991.4667  # CUL!J107 was 1687.8948
$ws.Cells.Item(107, 12).Value = 5974.4001  # CUL!L107 was 5063.6844
$ws.Cells.Item(107, 14).Value = -9814.400099999999  # CUL!N107 was -8903.6844

$ws.Cells.Item(137, 8).Value = 9039.857  # CUL!H137 was 9673.691999999999
$ws.Cells.Item(137, 10).Value = 2983.3333  # CUL!J137 was 3420
$ws.Cells.Item(137, 12).Value = 8949.999899999999  # CUL!L137 was 10260
$ws.Cells.Item(137, 14).Value = -19149.9999  # CUL!N137 was -20460

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0  # GSM!H5 was 1000
$ws.Cells.Item(5, 9).Value = 0  # GSM!I5 was 1000
$ws.Cells.Item(5, 11).Value = 0  # GSM!K5 was 1000
$ws.Cells.Item(5, 13).ClearContents()  # GSM!M5 was -888

$ws.Cells.Item(10, 8).Value = 250225  # GSM!H10 was 450
$ws.Cells.Item(10, 10).Value = 250225  # GSM!J10 was 450
$ws.Cells.Item(10, 12).Value = 250225  # GSM!L10 was 450
$ws.Cells.Item(10, 14).Value = -250563  # GSM!N10 was -788

$ws.Cells.Item(12, 8).Value = 500000  # GSM!H12 was 0
$ws.Cells.Item(12, 9).Value = 500000  # GSM!I12 was 0
$ws.Cells.Item(12, 11).Value = 500000  # GSM!K12 was 0
$ws.Cells.Item(12, 13).Value = -499860  # GSM!M12 was None

$ws.Cells.Item(14, 8).Value = 173121.8  # GSM!H14 was 250986.08
$ws.Cells.Item(14, 10).Value = 2648.077  # GSM!J14 was 1132.8572
$ws.Cells.Item(14, 12).Value = 2648.077  # GSM!L14 was 1132.8572
$ws.Cells.Item(14, 14).Value = -2984.077  # GSM!N14 was -1468.8572

$ws.Cells.Item(17, 8).Value = 6436.857  # GSM!H17 was 7373.9165
$ws.Cells.Item(17, 10).Value = 6662.769  # GSM!J17 was 7726.091
$ws.Cells.Item(17, 12).Value = 6662.769  # GSM!L17 was 7726.091
$ws.Cells.Item(17, 14).Value = -6998.769  # GSM!N17 was -8062.091

$ws.Cells.Item(19, 8).Value = 10979.8125  # GSM!H19 was 9881.091
$ws.Cells.Item(19, 10).Value = 12534.429  # GSM!J19 was 12055.223
$ws.Cells.Item(19, 12).Value = 12534.429  # GSM!L19 was 12055.223
$ws.Cells.Item(19, 14).Value = -13110.429  # GSM!N19 was -12631.223

$ws.Cells.Item(22, 8).Value = 5085.2856  # GSM!H22 was 5699.5
$ws.Cells.Item(22, 9).Value = 2799.6667  # GSM!I22 was 3166.3333
$ws.Cells.Item(22, 10).Value = 6799.5  # GSM!J22 was 8232.666999999999
$ws.Cells.Item(22, 11).Value = 2799.6667  # GSM!K22 was 3166.3333
$ws.Cells.Item(22, 12).Value = 6799.5  # GSM!L22 was 8232.666999999999
$ws.Cells.Item(22, 13).Value = -2270.6667  # GSM!M22 was -2637.3333
$ws.Cells.Item(22, 14).Value = -7857.5  # GSM!N22 was -9290.666999999999

$ws.Cells.Item(80, 8).Value = 6380.3335  # GSM!H80 was 5320.75
$ws.Cells.Item(80, 9).Value = 4999  # GSM!I80 was 3399.5
$ws.Cells.Item(80, 10).Value = 6656.6  # GSM!J80 was 5961.1665
$ws.Cells.Item(80, 11).Value = 4999  # GSM!K80 was 3399.5
$ws.Cells.Item(80, 12).Value = 6656.6  # GSM!L80 was 5961.1665
$ws.Cells.Item(80, 13).Value = -4001  # GSM!M80 was -2401.5
$ws.Cells.Item(80, 14).Value = -8652.6  # GSM!N80 was -7957.1665

$ws.Cells.Item(83, 8).Value = 6380.3335  # GSM!H83 was 5320.75
$ws.Cells.Item(83, 9).Value = 4999  # GSM!I83 was 3399.5
$ws.Cells.Item(83, 10).Value = 6656.6  # GSM!J83 was 5961.1665
$ws.Cells.Item(83, 11).Value = 24995  # GSM!K83 was 16997.5
$ws.Cells.Item(83, 12).Value = 33283  # GSM!L83 was 29805.8325
$ws.Cells.Item(83, 13).Value = -20003  # GSM!M83 was -12005.5
$ws.Cells.Item(83, 14).Value = -43267  # GSM!N83 was -39789.8325

$ws.Cells.Item(126, 8).Value = 3358.261  # GSM!H126 was 3265.348
$ws.Cells.Item(126, 9).Value = 2484.7273  # GSM!I126 was 2419.6
$ws.Cells.Item(126, 10).Value = 4159  # GSM!J126 was 3915.923
$ws.Cells.Item(126, 11).Value = 7454.1819  # GSM!K126 was 7258.799999999999
$ws.Cells.Item(126, 12).Value = 12477  # GSM!L126 was 11747.769
$ws.Cells.Item(126, 13).Value = -4984.1819  # GSM!M126 was -4788.799999999999
$ws.Cells.Item(126, 14).Value = -17417  # GSM!N126 was -16687.769

$ws.Cells.Item(132, 8).Value = 4885.579  # GSM!H132 was 5039.5557
$ws.Cells.Item(132, 9).Value = 3791.3  # GSM!I132 was 3898.3618
$ws.Cells.Item(132, 11).Value = 11373.9  # GSM!K132 was 11695.0854
$ws.Cells.Item(132, 13).Value = -8843.900000000001  # GSM!M132 was -9165.0854

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3799.6667  # LTW!H22 was 3902.4
$ws.Cells.Item(22, 9).Value = 2581.4614  # LTW!I22 was 2713.25
$ws.Cells.Item(22, 10).Value = 4488.2173  # LTW!J22 was 4522.826
$ws.Cells.Item(22, 11).Value = 2581.4614  # LTW!K22 was 2713.25
$ws.Cells.Item(22, 12).Value = 4488.2173  # LTW!L22 was 4522.826
$ws.Cells.Item(22, 13).Value = -2286.4614  # LTW!M22 was -2418.25
$ws.Cells.Item(22, 14).Value = -5078.2173  # LTW!N22 was -5112.826

$ws.Cells.Item(27, 8).Value = 3799.6667  # LTW!H27 was 3902.4
$ws.Cells.Item(27, 9).Value = 2581.4614  # LTW!I27 was 2713.25
$ws.Cells.Item(27, 10).Value = 4488.2173  # LTW!J27 was 4522.826
$ws.Cells.Item(27, 11).Value = 2581.4614  # LTW!K27 was 2713.25
$ws.Cells.Item(27, 12).Value = 4488.2173  # LTW!L27 was 4522.826
$ws.Cells.Item(27, 13).Value = -2474.4614  # LTW!M27 was -2606.25
$ws.Cells.Item(27, 14).Value = -4702.2173  # LTW!N27 was -4736.826

$ws.Cells.Item(40, 8).Value = 12605.223  # LTW!H40 was 13052.588
$ws.Cells.Item(40, 9).Value = 11806.934  # LTW!I40 was 12293.143
$ws.Cells.Item(40, 11).Value = 11806.934  # LTW!K40 was 12293.143
$ws.Cells.Item(40, 13).Value = -11670.934  # LTW!M40 was -12157.143

$ws.Cells.Item(82, 8).Value = 1049.7646  # LTW!H82 was 983.7727
$ws.Cells.Item(82, 9).Value = 1036.3636  # LTW!I82 was 989.1429000000001
$ws.Cells.Item(82, 10).Value = 1074.3334  # LTW!J82 was 974.375
$ws.Cells.Item(82, 11).Value = 1036.3636  # LTW!K82 was 989.1429000000001
$ws.Cells.Item(82, 12).Value = 1074.3334  # LTW!L82 was 974.375
$ws.Cells.Item(82, 13).Value = -675.3635999999999  # LTW!M82 was -628.1429000000001
$ws.Cells.Item(82, 14).Value = -1796.3334  # LTW!N82 was -1696.375

$ws.Cells.Item(85, 8).Value = 1049.7646  # LTW!H85 was 983.7727
$ws.Cells.Item(85, 9).Value = 1036.3636  # LTW!I85 was 989.1429000000001
$ws.Cells.Item(85, 10).Value = 1074.3334  # LTW!J85 was 974.375
$ws.Cells.Item(85, 11).Value = 1036.3636  # LTW!K85 was 989.1429000000001
$ws.Cells.Item(85, 12).Value = 1074.3334  # LTW!L85 was 974.375
$ws.Cells.Item(85, 13).Value = 211.6364000000001  # LTW!M85 was 258.8570999999999
$ws.Cells.Item(85, 14).Value = -3570.3334  # LTW!N85 was -3470.375

$ws.Cells.Item(136, 8).Value = 5272.08  # LTW!H136 was 5252.48
$ws.Cells.Item(136, 9).Value = 5041.1875  # LTW!I136 was 4909.294
$ws.Cells.Item(136, 10).Value = 5682.5557  # LTW!J136 was 5981.75
$ws.Cells.Item(136, 11).Value = 15123.5625  # LTW!K136 was 14727.882
$ws.Cells.Item(136, 12).Value = 17047.6671  # LTW!L136 was 17945.25
$ws.Cells.Item(136, 13).Value = -12573.5625  # LTW!M136 was -12177.882
$ws.Cells.Item(136, 14).Value = -22147.6671  # LTW!N136 was -23045.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4000  # WVR!H62 was 0
$ws.Cells.Item(62, 10).Value = 4000  # WVR!J62 was 0
$ws.Cells.Item(62, 12).Value = 4000  # WVR!L62 was 0
$ws.Cells.Item(62, 14).Value = -5248  # WVR!N62 was None

$ws.Cells.Item(65, 8).Value = 4000  # WVR!H65 was 0
$ws.Cells.Item(65, 10).Value = 4000  # WVR!J65 was 0
$ws.Cells.Item(65, 12).Value = 20000  # WVR!L65 was 0
$ws.Cells.Item(65, 14).Value = -26240  # WVR!N65 was None

$ws.Cells.Item(122, 8).Value = 1982.9474  # WVR!H122 was 1990.9524
$ws.Cells.Item(122, 9).Value = 1982.9474  # WVR!I122 was 2030.5
$ws.Cells.Item(122, 10).Value = 0  # WVR!J122 was 1200
$ws.Cells.Item(122, 11).Value = 5948.8422  # WVR!K122 was 6091.5
$ws.Cells.Item(122, 12).Value = 0  # WVR!L122 was 3600
$ws.Cells.Item(122, 13).Value = -3498.8422  # WVR!M122 was -3641.5
$ws.Cells.Item(122, 14).ClearContents()  # WVR!N122 was -8500

$ws.Cells.Item(132, 8).Value = 4580.2  # WVR!H132 was 4481.724
$ws.Cells.Item(132, 9).Value = 2846.682  # WVR!I132 was 2817.1365
$ws.Cells.Item(132, 10).Value = 9347.375  # WVR!J132 was 9713.286
$ws.Cells.Item(132, 11).Value = 8540.045999999998  # WVR!K132 was 8451.4095
$ws.Cells.Item(132, 12).Value = 28042.125  # WVR!L132 was 29139.858
$ws.Cells.Item(132, 13).Value = -6010.045999999998  # WVR!M132 was -5921.4095
$ws.Cells.Item(132, 14).Value = -33102.125  # WVR!N132 was -34199.858

$ws.Cells.Item(136, 8).Value = 3689.76  # WVR!H136 was 3688.8
$ws.Cells.Item(136, 9).Value = 1406.4166  # WVR!I136 was 1404.4166
$ws.Cells.Item(136, 11).Value = 4219.2498  # WVR!K136 was 4213.2498
$ws.Cells.Item(136, 13).Value = -1669.2498  # WVR!M136 was -1663.2498

